$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.124.71"
$ws.Range("E2").Value = "  +0.21%  "

$ws.Range("D3").Value = "'2.328.94"

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'304.95"
$ws.Range("E5").Value = "  +1.66%  "

$ws.Range("D6").Value = "'97.73"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("E9").Value = "  -0.56%  "

$ws.Range("D10").Value = "'35.61"
$ws.Range("E10").Value = "  -0.68%  "

$ws.Range("D11").Value = "'19.54"
$ws.Range("E11").Value = "  +7.59%  "

$ws.Range("D12").Value = "'0.0797"
$ws.Range("E12").Value = "  +0.95%  "

$ws.Range("E13").Value = "  +0.82%  "

$ws.Range("E14").Value = "  +1.84%  "

$ws.Range("D15").Value = "'2.692.16"
$ws.Range("E15").Value = "  +1.05%  "

$ws.Range("D16").Value = "'2.318.63"
$ws.Range("E16").Value = "  +0.95%  "

$ws.Range("E17").Value = "  +0.99%  "

$ws.Range("D18").Value = "'43.024.15"
$ws.Range("E18").Value = "  +0.14%  "

$ws.Range("D19").Value = "'12.57"
$ws.Range("E19").Value = "  -1.21%  "

$ws.Range("E20").Value = "  -0.40%  "

$ws.Range("E21").Value = "  +0.45%  "

$ws.Range("D22").Value = "'67.97"
$ws.Range("E22").Value = "  -0.14%  "

$ws.Range("D23").Value = "'237.62"
$ws.Range("E23").Value = "  -1.28%  "

$ws.Range("E24").Value = "  +3.11%  "

$ws.Range("E25").Value = "  +0.77%  "

$ws.Range("E26").Value = "  +0.19%  "

$ws.Range("D27").Value = "'24.96"
$ws.Range("E27").Value = "  -2.24%  "

$ws.Range("D28").Value = "'166.23"
$ws.Range("E28").Value = "  +0.24%  "

$ws.Range("E29").Value = "  +1.98%  "

$ws.Range("E30").Value = "  +0.77%  "

$ws.Range("D31").Value = "'33.23"
$ws.Range("E31").Value = "  +0.12%  "

$ws.Range("D33").Value = "'18.04"
$ws.Range("E33").Value = "  +5.78%  "

$ws.Range("E34").Value = "  -0.72%  "

$ws.Range("E35").Value = "  -8.47%  "

$ws.Range("D37").Value = "'0.0697"
$ws.Range("E37").Value = "  +1.54%  "

$ws.Range("E38").Value = "  +0.21%  "

$ws.Range("D39").Value = "'2.80"
$ws.Range("E39").Value = "  +2.11%  "

$ws.Range("D40").Value = "'1.77"
$ws.Range("E40").Value = "  -0.08%  "

$ws.Range("D41").Value = "'0.110"
$ws.Range("E41").Value = "  -0.35%  "

$ws.Range("D42").Value = "'1.994.99"
$ws.Range("E42").Value = "  -0.80%  "

$ws.Range("D43").Value = "'10.72"
$ws.Range("E43").Value = "  +5.71%  "

$ws.Range("E44").Value = "  -0.16%  "

$ws.Range("D45").Value = "'18.15"
$ws.Range("E45").Value = "  +4.53%  "

$ws.Range("D46").Value = "'2.08"
$ws.Range("E46").Value = "  -3.93%  "

$ws.Range("E47").Value = "  -0.71%  "

$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").Value = "'2.89"
$ws.Range("E48").Value = "  +2.25%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "'2.558.40"
$ws.Range("E49").Value = "  +1.02%  "

$ws.Range("D50").Value = "'53.73"
$ws.Range("E50").Value = "  -0.14%  "

$ws.Range("D51").Value = "'72.01"
$ws.Range("E51").Value = "  -0.64%  "
